$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers for Wins, Losses, Ties in columns AD, AE, AF (30, 31, 32)
$ws.Cells.Item(1, 30).Value = "Wins"
$ws.Cells.Item(1, 31).Value = "Losses"
$ws.Cells.Item(1, 32).Value = "Ties"

# Copy the style of an existing header cell (AC1) to the new header cells
$ws.Range("AC1").Copy()
$headerRange = $ws.Range($ws.Cells.Item(1, 30), $ws.Cells.Item(1, 32))
$headerRange.PasteSpecial(-4122)  # xlPasteFormats

# Fill in the Wins/Losses/Ties values for each data row (2 through 61)
$lastRow = 61
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 30).Value = 75
    $ws.Cells.Item($r, 31).Value = 87
    $ws.Cells.Item($r, 32).Value = 0
}
